# Insert a new record row at row 33 (shifting all existing rows 33..105 down to 34..106)
# and populate it with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 33; this pushes rows 33-105 down to 34-106
# and extends the sheet dimension to A1:T106.
$ws.Rows.Item(33).Insert()

# Populate the new row 33 with the new record's values.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,T keep the same values as the (now shifted) row 34,
# since the Insert operation does not copy values into the newly blank row.
$ws.Range("A33").Value = 10
$ws.Range("B33").Value = "Vega Modelo de Temuco"
$ws.Range("C33").Value = "La Araucanía"
$ws.Range("D33").Value = "7/18/2023"
$ws.Range("E33").Value = 9
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100108
$ws.Range("H33").Value = "Tropicales y subtropicales"
$ws.Range("I33").Value = 100108003
$ws.Range("J33").Value = "Maracuyá"
$ws.Range("K33").Value = "Sin especificar"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 15
$ws.Range("N33").Value = 45000
$ws.Range("O33").Value = 45000
$ws.Range("P33").Value = 45000
$ws.Range("Q33").Value = '$/caja 18 kilos'
$ws.Range("R33").Value = "Región de Arica y Parinacota"
$ws.Range("S33").Value = 2500
$ws.Range("T33").Value = 18
